# Add an additional intermediate acid dose row to the "Slurry" sheet.
# 5.7 kg/t acid dose -> 25% reduction in pH relative to un-treated digestate.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Slurry")

# New row 5, same shape/format as rows 2-4 (copy formatting from row 4).
$ws.Range("A4:E4").Copy($ws.Range("A5:E5")) | Out-Null

$ws.Range("A5").Value = "Afgasset biomasse"
$ws.Range("B5").Value = "Digestate"
$ws.Range("C5").Value = "5.7 kg/t"
$ws.Range("D5").Value = 5.1
$ws.Range("E5").Formula = "=7.9-1.01"

# Match the new active cell/selection recorded in the saved file.
$ws.Range("E5").Select() | Out-Null
